# Append the new Nifty 500 daily data rows (2024-09-30 .. 2025-01-14)
# to Sheet1, rows 3408-3480, columns A:E, matching the source row layout:
# A = Date (serial, formatted as "yyyy-mm-dd h:mm:ss"), B = Total Returns Index,
# C = P/E, D = P/B, E = Div Yield.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3408

$newRows = @(
    @(45565,38261.39,27.87,4.22,1.04),
    @(45566,38302.54,27.9,4.22,1.04),
    @(45568,37488.51,27.31,4.13,1.06),
    @(45569,37140.49,27.05,4.09,1.07),
    @(45572,36583.47,26.65,4.03,1.09),
    @(45573,37105.17,27.03,4.09,1.07),
    @(45574,37242.54,27.13,4.11,1.07),
    @(45575,37240.29,27.12,4.11,1.07),
    @(45576,37261.17,27.14,4.11,1.07),
    @(45579,37443.97,27.27,4.13,1.06),
    @(45580,37457.6,27.27,4.13,1.06),
    @(45581,37358.9,27.2,4.12,1.06),
    @(45582,36891.53,26.86,4.07,1.08),
    @(45583,36994.04,26.94,4.08,1.08),
    @(45586,36713.97,26.68,4.05,1.08),
    @(45587,36008.15,26.07,3.97,1.1),
    @(45588,36050.23,26.1,3.98,1.1),
    @(45589,35966.15,26,3.97,1.1),
    @(45590,35510.88,25.65,3.92,1.08),
    @(45593,35749.31,25.85,3.94,1.08),
    @(45594,35992.9,26.26,3.97,1.11),
    @(45595,35921.77,26.21,3.96,1.11),
    @(45596,35830.98,26.07,3.95,1.13),
    @(45597,36042.89,26.23,3.97,1.12),
    @(45600,35561.76,25.85,3.92,1.12),
    @(45601,35816.9,26.03,3.95,1.13),
    @(45602,36377.55,26.43,4.01,1.12),
    @(45603,36023.67,26.16,3.97,1.13),
    @(45604,35771.65,25.8,3.94,1.15),
    @(45607,35685.95,25.71,3.93,1.14),
    @(45608,35275.39,25.26,3.89,1.15),
    @(45609,34618.08,24.83,3.81,1.17),
    @(45610,34699.38,24.88,3.82,1.18),
    @(45614,34599.09,24.78,3.81,1.17),
    @(45615,34766.49,24.87,3.83,1.18),
    @(45617,34476.59,24.78,3.8,1.16),
    @(45618,35116,25.24,3.87,1.14),
    @(45621,35635.18,25.61,3.92,1.12),
    @(45622,35649.51,25.62,3.92,1.13),
    @(45623,35853.83,25.77,3.93,1.12),
    @(45624,35573.8,25.57,3.89,1.13),
    @(45625,35847.4,25.77,3.92,1.12),
    @(45628,36082.61,25.94,3.95,1.11),
    @(45629,36374.88,26.14,3.98,1.11),
    @(45630,36499.26,26.23,3.99,1.1),
    @(45631,36781.26,26.44,4.02,1.09),
    @(45632,36835.17,26.48,4.03,1.09),
    @(45635,36803.45,26.45,4.02,1.09),
    @(45636,36839.6,26.48,4.03,1.09),
    @(45637,36907.95,26.53,4.04,1.09),
    @(45638,36747.2,26.41,4.02,1.08),
    @(45639,36909.23,26.53,4.04,1.08),
    @(45642,36898.98,26.52,4.03,1.08),
    @(45643,36493.34,26.23,3.99,1.09),
    @(45644,36238.39,26.05,3.96,1.1),
    @(45645,35950.45,25.84,3.93,1.11),
    @(45646,35266.62,25.35,3.86,1.13),
    @(45649,35413.41,25.45,3.87,1.12),
    @(45650,35401.02,25.44,3.86,1.13),
    @(45652,35444.76,25.47,3.87,1.13),
    @(45653,35468.25,25.49,3.87,1.12),
    @(45656,35329.04,25.39,3.85,1.12),
    @(45657,35357.9,25.56,3.88,1.12),
    @(45658,35526.08,25.68,3.9,1.11),
    @(45659,36060.05,26.07,3.96,1.1),
    @(45660,35883.72,25.94,3.94,1.1),
    @(45663,35120.87,25.39,3.85,1.13),
    @(45664,35304.55,25.52,3.87,1.12),
    @(45665,35131.01,25.4,3.85,1.12),
    @(45666,34838.29,25.19,3.82,1.13),
    @(45667,34447.06,24.9,3.78,1.15),
    @(45670,33569.22,24.26,3.68,1.18),
    @(45671,33944.68,24.53,3.72,1.16)
)

$nRows = $newRows.Count
$nCols = 5
$endRow = $startRow + $nRows - 1

$arr = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i,$j] = $newRows[$i][$j]
    }
}

$ws.Range("A${startRow}:E${endRow}").Value = $arr

# Column A keeps the same date/time number format used by the preceding
# rows in the sheet (style index 5 -> numFmtId 167, "yyyy-mm-dd h:mm:ss").
$ws.Range("A${startRow}:A${endRow}").NumberFormat = "yyyy-mm-dd h:mm:ss"
